# Apply updated market/profit calculations to each job sheet
# (values refreshed by the scheduled market-data runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1302.7273
$ws.Range("I32").Value = 1291.625
$ws.Range("J32").Value = 1332.3334
$ws.Range("K32").Value = 1291.625
$ws.Range("L32").Value = 1332.3334
$ws.Range("M32").Value = -965.625
$ws.Range("N32").Value = -1984.3334
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 2525
$ws.Range("I100").Value = 2030.3572
$ws.Range("K100").Value = 2030.3572
$ws.Range("M100").Value = -1489.3572
$ws.Range("H103").Value = 593.3333
$ws.Range("J103").Value = 593.3333
$ws.Range("L103").Value = 1779.9999
$ws.Range("N103").Value = -2951.9999
$ws.Range("H107").Value = 620
$ws.Range("J107").Value = 375
$ws.Range("L107").Value = 375
$ws.Range("N107").Value = -4215
$ws.Range("H112").Value = 2836.5
$ws.Range("J112").Value = 2836.5
$ws.Range("L112").Value = 8509.5
$ws.Range("N112").Value = -10725.5
$ws.Range("H116").Value = 5399.7
$ws.Range("I116").Value = 2999.625
$ws.Range("K116").Value = 2999.625
$ws.Range("M116").Value = 442.375
$ws.Range("H137").Value = 10754795
$ws.Range("I137").Value = 16668400
$ws.Range("J137").Value = 2786.3635
$ws.Range("K137").Value = 50005200
$ws.Range("L137").Value = 8359.0905
$ws.Range("M137").Value = -50002650
$ws.Range("N137").Value = -13459.0905
$ws.Range("H138").Value = 6448
$ws.Range("I138").Value = 6601.8184
$ws.Range("J138").Value = 6387.5713
$ws.Range("K138").Value = 19805.4552
$ws.Range("L138").Value = 19162.7139
$ws.Range("M138").Value = -14665.4552
$ws.Range("N138").Value = -29442.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 535.44446
$ws.Range("I5").Value = 331.2857
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 331.2857
$ws.Range("L5").Value = 1250
$ws.Range("M5").Value = -219.2857
$ws.Range("N5").Value = -1474
$ws.Range("H32").Value = 14258.6
$ws.Range("I32").Value = 11067.692
$ws.Range("K32").Value = 11067.692
$ws.Range("M32").Value = -10780.692
$ws.Range("H74").Value = 4061957.5
$ws.Range("I74").Value = 4931306
$ws.Range("K74").Value = 4931306
$ws.Range("M74").Value = -4930432
$ws.Range("H77").Value = 4061957.5
$ws.Range("I77").Value = 4931306
$ws.Range("K77").Value = 24656530
$ws.Range("M77").Value = -24652162
$ws.Range("H88").Value = 980.75
$ws.Range("I88").Value = 766.6667
$ws.Range("J88").Value = 1109.2
$ws.Range("K88").Value = 766.6667
$ws.Range("L88").Value = 1109.2
$ws.Range("M88").Value = -360.6667
$ws.Range("N88").Value = -1921.2
$ws.Range("H91").Value = 980.75
$ws.Range("I91").Value = 766.6667
$ws.Range("J91").Value = 1109.2
$ws.Range("K91").Value = 766.6667
$ws.Range("L91").Value = 1109.2
$ws.Range("M91").Value = 637.3333
$ws.Range("N91").Value = -3917.2
$ws.Range("H96").Value = 19661.572
$ws.Range("J96").Value = 19661.572
$ws.Range("L96").Value = 19661.572
$ws.Range("N96").Value = -25153.572
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 535.44446
$ws.Range("I4").Value = 331.2857
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 331.2857
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -216.2857
$ws.Range("N4").Value = -1480
$ws.Range("H75").Value = 8000
$ws.Range("I75").Value = 8000
$ws.Range("K75").Value = 8000
$ws.Range("M75").Value = -7064
$ws.Range("H78").Value = 8000
$ws.Range("I78").Value = 8000
$ws.Range("K78").Value = 24000
$ws.Range("M78").Value = -19320
$ws.Range("H86").Value = 19670804
$ws.Range("I86").Value = 70909.234
$ws.Range("K86").Value = 70909.234
$ws.Range("M86").Value = -69786.234
$ws.Range("H89").Value = 19670804
$ws.Range("I89").Value = 70909.234
$ws.Range("K89").Value = 354546.17
$ws.Range("M89").Value = -348930.17
$ws.Range("H134").Value = 1066.1111
$ws.Range("I134").Value = 949.375
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2848.125
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -313.125
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 58823804
$ws.Range("I7").Value = 76923304
$ws.Range("K7").Value = 76923304
$ws.Range("M7").Value = -76923191
$ws.Range("H31").Value = 3312.5
$ws.Range("I31").Value = 3312.5
$ws.Range("K31").Value = 3312.5
$ws.Range("M31").Value = -3017.5
$ws.Range("H34").Value = 3312.5
$ws.Range("I34").Value = 3312.5
$ws.Range("K34").Value = 3312.5
$ws.Range("M34").Value = -3110.5
$ws.Range("H47").Value = 24999.75
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 24999.75
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 24999.75
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -26131.75
$ws.Range("H62").Value = 3570.5
$ws.Range("I62").Value = 2899
$ws.Range("J62").Value = 5585
$ws.Range("K62").Value = 2899
$ws.Range("L62").Value = 5585
$ws.Range("M62").Value = -2275
$ws.Range("N62").Value = -6833
$ws.Range("H65").Value = 3570.5
$ws.Range("I65").Value = 2899
$ws.Range("J65").Value = 5585
$ws.Range("K65").Value = 14495
$ws.Range("L65").Value = 27925
$ws.Range("M65").Value = -11375
$ws.Range("N65").Value = -34165
$ws.Range("H86").Value = 9094.666999999999
$ws.Range("I86").Value = 9094.666999999999
$ws.Range("K86").Value = 9094.666999999999
$ws.Range("M86").Value = -7971.666999999999
$ws.Range("H89").Value = 9094.666999999999
$ws.Range("I89").Value = 9094.666999999999
$ws.Range("K89").Value = 45473.335
$ws.Range("M89").Value = -39857.335
$ws.Range("H107").Value = 900.53845
$ws.Range("J107").Value = 1367
$ws.Range("L107").Value = 1367
$ws.Range("N107").Value = -5207
$ws.Range("H141").Value = 627536.9
$ws.Range("J141").Value = 627536.9
$ws.Range("L141").Value = 627536.9
$ws.Range("N141").Value = -637896.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1075.6666
$ws.Range("I109").Value = 1075.6666
$ws.Range("K109").Value = 3226.9998
$ws.Range("M109").Value = -2186.9998
$ws.Range("H113").Value = 3768.0833
$ws.Range("I113").Value = 3915.3333
$ws.Range("J113").Value = 3719
$ws.Range("K113").Value = 11745.9999
$ws.Range("L113").Value = 11157
$ws.Range("M113").Value = -9575.999899999999
$ws.Range("N113").Value = -15497
$ws.Range("H131").Value = 1500
$ws.Range("I131").Value = 1500
$ws.Range("K131").Value = 4500
$ws.Range("M131").Value = 540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 22000
$ws.Range("I22").Value = 18000
$ws.Range("K22").Value = 18000
$ws.Range("M22").Value = -17471
$ws.Range("H95").Value = 19097
$ws.Range("J95").Value = 19097
$ws.Range("L95").Value = 19097
$ws.Range("N95").Value = -24589
$ws.Range("H132").Value = 22226018
$ws.Range("I132").Value = 4065.6428
$ws.Range("K132").Value = 12196.9284
$ws.Range("M132").Value = -9666.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2499
$ws.Range("I16").Value = 2499
$ws.Range("K16").Value = 2499
$ws.Range("M16").Value = -2329
$ws.Range("H22").Value = 33335198
$ws.Range("I22").Value = 593.1667
$ws.Range("J22").Value = 55558268
$ws.Range("K22").Value = 593.1667
$ws.Range("L22").Value = 55558268
$ws.Range("M22").Value = -298.1667
$ws.Range("N22").Value = -55558858
$ws.Range("H27").Value = 33335198
$ws.Range("I27").Value = 593.1667
$ws.Range("J27").Value = 55558268
$ws.Range("K27").Value = 593.1667
$ws.Range("L27").Value = 55558268
$ws.Range("M27").Value = -486.1667
$ws.Range("N27").Value = -55558482
$ws.Range("H100").Value = 2998
$ws.Range("I100").Value = 2747.75
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 2747.75
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -2206.75
$ws.Range("N100").Value = -6082
$ws.Range("H101").Value = 22928
$ws.Range("J101").Value = 22928
$ws.Range("L101").Value = 22928
$ws.Range("N101").Value = -29418
$ws.Range("H132").Value = 2572.125
$ws.Range("I132").Value = 2572.125
$ws.Range("K132").Value = 7716.375
$ws.Range("M132").Value = -5186.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 8499
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 8499
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 8499
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -8779
$ws.Range("H54").Value = 4022079
$ws.Range("I54").Value = 4516089
$ws.Range("K54").Value = 4516089
$ws.Range("M54").Value = -4515569
$ws.Range("H132").Value = 58824400
$ws.Range("I132").Value = 924.8125
$ws.Range("K132").Value = 2774.4375
$ws.Range("M132").Value = -244.4375
